# Update country order (Pais column) and updated case numbers
# Spain covid-19 "paises.xlsx" daily refresh: 30 Marzo 2020 15:20 -> 15:50
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp banner in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 15:50"

# --- Re-order country names (column A) for rows 180-205 ---
# (the daily source re-sorted the country list; the row data stayed put
#  but the country label shown in column A for each row changed)
$ws.Range("A181").Value = "Santa Sede"
$ws.Range("A183").Value = "Benin"
$ws.Range("A186").Value = "San Bartolome"
$ws.Range("A187").Value = "Montserrat"
$ws.Range("A188").Value = "Fiyi"
$ws.Range("A189").Value = "Republica del Chad"
$ws.Range("A190").Value = "Nepal"
$ws.Range("A191").Value = "Mauritania"
$ws.Range("A192").Value = "Butan"
$ws.Range("A193").Value = "Islas Turcas y Caicos"
$ws.Range("A195").Value = "Nicaragua"
$ws.Range("A196").Value = "Republica de Africa Central"
$ws.Range("A198").Value = "Liberia"
$ws.Range("A199").Value = "Belice"
$ws.Range("A202").Value = "Islas Virgenes Britanicas"

# --- Updated case figures (Casos totales/Nuevos casos/Casos activos/Recuperados/Casos criticos/Muertes hoy/Muertes) ---
# Row 22
$ws.Range("B22").Value = 4330
$ws.Range("C22").Value = 74
$ws.Range("E22").Value = 4070
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 140
# Row 29
$ws.Range("B29").Value = 2449
$ws.Range("C29").Value = 310
$ws.Range("D29").Value = 156
$ws.Range("E29").Value = 2285
$ws.Range("F29").Value = 14
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 8
# Row 39
$ws.Range("D39").Value = 115
$ws.Range("E39").Value = 1330
# Row 51
$ws.Range("D51").Value = 228
$ws.Range("E51").Value = 570
# Row 53
$ws.Range("F53").Value = 62
# Row 60
$ws.Range("D60").Value = 61
$ws.Range("E60").Value = 545
# Row 81
$ws.Range("D81").Value = 21
$ws.Range("E81").Value = 280
# Row 106
$ws.Range("E106").Value = 105
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 2
# Row 113
$ws.Range("B113").Value = 100
$ws.Range("C113").Value = 9
$ws.Range("E113").Value = 82
# Row 114
$ws.Range("E114").Value = 92
$ws.Range("G114").Value = 3
$ws.Range("H114").Value = 4
# Row 141
$ws.Range("F141").Value = 5
# Row 189
$ws.Range("C189").Value = 2
$ws.Range("D189").Value = 0
$ws.Range("E189").Value = 5
# Row 190
$ws.Range("D190").Value = 1
$ws.Range("E190").Value = 4
# Row 191
$ws.Range("B191").Value = 5
$ws.Range("D191").Value = 2
$ws.Range("E191").Value = 3
# Row 193
$ws.Range("E193").Value = 4
$ws.Range("H193").Value = 0
# Row 195
$ws.Range("B195").Value = 4
$ws.Range("H195").Value = 1

